$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 with the same formatting (bold, border) as the other headers
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

# Data cells F2:F9 with the time_taken values
$times = @(
    "2021-10-05 13:40:59.552226",
    "2021-10-05 13:40:59.552238",
    "2021-10-05 13:40:59.552241",
    "2021-10-05 13:40:59.552244",
    "2021-10-05 13:40:59.552247",
    "2021-10-05 13:40:59.552250",
    "2021-10-05 13:40:59.552253",
    "2021-10-05 13:40:59.552256"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
